# Apply scheduled profit-data refresh for Midgardsormr sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 135
$ws.Range("H135").Value = 14584.919
$ws.Range("I135").Value = 503.51614
$ws.Range("K135").Value = 4531.64526
$ws.Range("M135").Value = -1996.64526

# ALC row 137
$ws.Range("H137").Value = 14217.357
$ws.Range("I137").Value = 16299.772
$ws.Range("K137").Value = 48899.31600000001
$ws.Range("M137").Value = -46349.31600000001

# ALC row 138
$ws.Range("H138").Value = 20507.414
$ws.Range("I138").Value = 2329.0334
$ws.Range("J138").Value = 39984.25
$ws.Range("K138").Value = 6987.100199999999
$ws.Range("L138").Value = 119952.75
$ws.Range("M138").Value = -1847.100199999999
$ws.Range("N138").Value = -130232.75

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 31243.982
$ws.Range("I32").Value = 23551.29
$ws.Range("K32").Value = 23551.29
$ws.Range("M32").Value = -23264.29

# ARM row 45
$ws.Range("H45").Value = 4402.4116
$ws.Range("I45").Value = 3038.5557
$ws.Range("J45").Value = 5936.75
$ws.Range("K45").Value = 3038.5557
$ws.Range("L45").Value = 5936.75
$ws.Range("M45").Value = -2661.5557
$ws.Range("N45").Value = -6690.75

# ARM row 61
$ws.Range("H61").Value = 4118.325
$ws.Range("I61").Value = 1230.7188
$ws.Range("J61").Value = 15668.75
$ws.Range("K61").Value = 1230.7188
$ws.Range("L61").Value = 15668.75
$ws.Range("M61").Value = -1018.7188
$ws.Range("N61").Value = -16092.75

# ARM row 74
$ws.Range("H74").Value = 127823.25
$ws.Range("I74").Value = 150728.22
$ws.Range("J74").Value = 13298.375
$ws.Range("K74").Value = 150728.22
$ws.Range("L74").Value = 13298.375
$ws.Range("M74").Value = -149854.22
$ws.Range("N74").Value = -15046.375

# ARM row 77
$ws.Range("H77").Value = 127823.25
$ws.Range("I77").Value = 150728.22
$ws.Range("J77").Value = 13298.375
$ws.Range("K77").Value = 753641.1
$ws.Range("L77").Value = 66491.875
$ws.Range("M77").Value = -749273.1
$ws.Range("N77").Value = -75227.875

# ARM row 132
$ws.Range("H132").Value = 1459.9385
$ws.Range("I132").Value = 1173.8684
$ws.Range("J132").Value = 1862.5555
$ws.Range("K132").Value = 3521.6052
$ws.Range("L132").Value = 5587.666499999999
$ws.Range("M132").Value = -991.6052
$ws.Range("N132").Value = -10647.6665

# ARM row 136
$ws.Range("H136").Value = 4118.325
$ws.Range("I136").Value = 1230.7188
$ws.Range("J136").Value = 15668.75
$ws.Range("K136").Value = 3692.1564
$ws.Range("L136").Value = 47006.25
$ws.Range("M136").Value = -1142.1564
$ws.Range("N136").Value = -52106.25

# ARM row 141
$ws.Range("H141").Value = 120000
$ws.Range("J141").Value = 120000
$ws.Range("L141").Value = 120000
$ws.Range("N141").Value = -130360

$ws = $wb.Worksheets.Item("BSM")
# BSM row 94
$ws.Range("H94").Value = 5520.696
$ws.Range("I94").Value = 7158.067
$ws.Range("K94").Value = 7158.067
$ws.Range("M94").Value = -6707.067

# BSM row 134
$ws.Range("H134").Value = 2063.4375
$ws.Range("I134").Value = 1223.3478
$ws.Range("J134").Value = 4210.3335
$ws.Range("K134").Value = 3670.0434
$ws.Range("L134").Value = 12631.0005
$ws.Range("M134").Value = -1135.0434
$ws.Range("N134").Value = -17701.0005

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 3127794.5
$ws.Range("I31").Value = 4762977
$ws.Range("K31").Value = 4762977
$ws.Range("M31").Value = -4762682

# CRP row 34
$ws.Range("H34").Value = 3127794.5
$ws.Range("I34").Value = 4762977
$ws.Range("K34").Value = 4762977
$ws.Range("M34").Value = -4762775

# CRP row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

# CRP row 105
$ws.Range("H105").Value = 1492.5714
$ws.Range("I105").Value = 1034.25
$ws.Range("K105").Value = 1034.25
$ws.Range("M105").Value = 712.75

# CRP row 132
$ws.Range("H132").Value = 32520.625
$ws.Range("I132").Value = 36745
$ws.Range("K132").Value = 110235
$ws.Range("M132").Value = -107705

# CRP row 134
$ws.Range("H134").Value = 1610.0869
$ws.Range("I134").Value = 1533.3684
$ws.Range("K134").Value = 4600.1052
$ws.Range("M134").Value = -2065.1052

$ws = $wb.Worksheets.Item("CUL")
# CUL row 107
$ws.Range("H107").Value = 3428.6365
$ws.Range("I107").Value = 7828.75
$ws.Range("J107").Value = 914.2857
$ws.Range("K107").Value = 23486.25
$ws.Range("L107").Value = 2742.8571
$ws.Range("M107").Value = -21566.25
$ws.Range("N107").Value = -6582.8571

# CUL row 131
$ws.Range("H131").Value = 2000.075
$ws.Range("J131").Value = 2399.3333
$ws.Range("L131").Value = 7197.999899999999
$ws.Range("N131").Value = -17277.9999

$ws = $wb.Worksheets.Item("GSM")
# GSM row 132
$ws.Range("H132").Value = 2182.6182
$ws.Range("I132").Value = 2222
$ws.Range("K132").Value = 6666
$ws.Range("M132").Value = -4136

$ws = $wb.Worksheets.Item("LTW")
# LTW row 46
$ws.Range("H46").Value = 3640.4443
$ws.Range("I46").Value = 800
$ws.Range("K46").Value = 800
$ws.Range("M46").Value = -612

# LTW row 51
$ws.Range("H51").Value = 140000
$ws.Range("J51").Value = 140000
$ws.Range("L51").Value = 140000
$ws.Range("N51").Value = -140956

# LTW row 55
$ws.Range("H55").Value = 676.625
$ws.Range("I55").Value = 302.875
$ws.Range("J55").Value = 1050.375
$ws.Range("K55").Value = 302.875
$ws.Range("L55").Value = 1050.375
$ws.Range("M55").Value = -129.875
$ws.Range("N55").Value = -1396.375

# LTW row 132
$ws.Range("H132").Value = 1655.742
$ws.Range("I132").Value = 1199.7593
$ws.Range("J132").Value = 4733.625
$ws.Range("K132").Value = 3599.2779
$ws.Range("L132").Value = 14200.875
$ws.Range("M132").Value = -1069.2779
$ws.Range("N132").Value = -19260.875

# LTW row 136
$ws.Range("H136").Value = 3050.762
$ws.Range("I136").Value = 3065.7368
$ws.Range("J136").Value = 2908.5
$ws.Range("K136").Value = 9197.2104
$ws.Range("L136").Value = 8725.5
$ws.Range("M136").Value = -6647.2104
$ws.Range("N136").Value = -13825.5

$ws = $wb.Worksheets.Item("WVR")
# WVR row 2
$ws.Range("H2").Value = 2714.375
$ws.Range("I2").Value = 2714.375
$ws.Range("K2").Value = 2714.375
$ws.Range("M2").Value = -2602.375

# WVR row 132
$ws.Range("H132").Value = 6965974.5
$ws.Range("I132").Value = 10447418
$ws.Range("K132").Value = 31342254
$ws.Range("M132").Value = -31339724

# WVR row 136
$ws.Range("H136").Value = 13457.615
$ws.Range("I136").Value = 14363.766
$ws.Range("J136").Value = 4939.8
$ws.Range("K136").Value = 43091.298
$ws.Range("L136").Value = 14819.4
$ws.Range("M136").Value = -40541.298
$ws.Range("N136").Value = -19919.4
